# Add a new "OPRF_PRINT_VAL // test the OPRF" line right after the existing
# "UNITTEST_ROUND2 ..." flag-description line in the preformatted flag list.
$d = $word.ActiveDocument

# Locate the paragraph that holds the UNITTEST_ROUND2 line.
$rng = $d.Content
$found = $rng.Find.Execute("UNITTEST_ROUND2", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$flagPara = $rng.Paragraphs(1)

# Insert a brand-new paragraph right after it (Word clones the paragraph /
# run formatting of $flagPara, matching the style used throughout this list).
$flagPara.Range.InsertParagraphAfter()
$newPara = $flagPara.Next()

# Populate the new paragraph with two runs - "OPRF_PRINT_VAL" and
# " // test the OPRF" - using the same JetBrains Mono / bold / green
# formatting as the rest of the flag list. Driving this through WordOpenXML
# keeps the two pieces of text as distinct runs (matching how the document
# was actually authored) instead of Word silently coalescing them into one.
$xmlFrag = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
    '<w:p>' + `
      '<w:pPr>' + `
        '<w:pStyle w:val="HTMLPreformatted"/>' + `
        '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' + `
        '<w:rPr>' + `
          '<w:rFonts w:ascii="JetBrains Mono" w:hAnsi="JetBrains Mono"/>' + `
          '<w:b/>' + `
          '<w:bCs/>' + `
          '<w:color w:val="1F542E"/>' + `
        '</w:rPr>' + `
      '</w:pPr>' + `
      '<w:r>' + `
        '<w:rPr>' + `
          '<w:rFonts w:ascii="JetBrains Mono" w:hAnsi="JetBrains Mono"/>' + `
          '<w:b/>' + `
          '<w:bCs/>' + `
          '<w:color w:val="1F542E"/>' + `
        '</w:rPr>' + `
        '<w:t>OPRF_PRINT_VAL</w:t>' + `
      '</w:r>' + `
      '<w:r>' + `
        '<w:rPr>' + `
          '<w:rFonts w:ascii="JetBrains Mono" w:hAnsi="JetBrains Mono"/>' + `
          '<w:b/>' + `
          '<w:bCs/>' + `
          '<w:color w:val="1F542E"/>' + `
        '</w:rPr>' + `
        '<w:t xml:space="preserve"> // test the OPRF</w:t>' + `
      '</w:r>' + `
    '</w:p>' + `
  '</w:body>' + `
'</w:document>'

[void]$newPara.Range.InsertXML($xmlFrag)
